$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 2, pushing the existing data rows down by one.
$ws.Rows("2:2").Insert()

# The inserted row copies formatting from the header row above it; reset it so it
# matches the plain (unstyled) look of the other data rows.
$ws.Range("A2:R2").ClearFormats()

# Populate the new weekly data record in row 2.
$ws.Cells.Item(2, 1).Value = 7
$ws.Cells.Item(2, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(2, 3).Value = "Ñuble"
$ws.Cells.Item(2, 4).Value = 44631
$ws.Cells.Item(2, 5).Value = 16
$ws.Cells.Item(2, 6).Value = 100112040
$ws.Cells.Item(2, 7).Value = "Cilantro"
$ws.Cells.Item(2, 8).Value = "Sin especificar"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 200
$ws.Cells.Item(2, 11).Value = 550
$ws.Cells.Item(2, 12).Value = 600
$ws.Cells.Item(2, 13).Value = 575
$ws.Cells.Item(2, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(2, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(2, 16).Value = 575
$ws.Cells.Item(2, 17).Value = 1
$ws.Cells.Item(2, 18).Value = "Hortaliza"

# Column D (Fecha) carries the same custom date/time number format as the rest
# of the date column; copy it from the row below (now row 3) onto the new cell.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat
